# Update plots for each sample
#
# The wildtype peak height for marker CYP2D6_10B (CYP2D6_002) in sample S1
# was re-measured (700 instead of 1000). With the corrected height the peak
# is now detected, which flips the called genotype for that marker from
# homozygous-mutant (TT) to heterozygous (CT), and the overall sample
# genotype from *10B/*10B to *1/*10B.

$wb = $excel.ActiveWorkbook

# --- peak_table: wildtype peak height for CYP2D6_10B / S1 -------------------
$peakTable = $wb.Worksheets.Item("peak_table")
$peakTable.Range("N3").Value = 700

# --- allele_table: recomputed detection results for the same peak -----------
$alleleTable = $wb.Worksheets.Item("allele_table")
$alleleTable.Range("K4").Value = 700
$alleleTable.Range("L4").Value = 1
$alleleTable.Range("M4").Value = $true
$alleleTable.Range("N4").Value = 40
$alleleTable.Range("O4").Value = 32.59
$alleleTable.Range("P4").Value = 748
$alleleTable.Range("Q4").Value = "ok"
$alleleTable.Range("R4").Value = ""

# --- marker_table: genotype/phenotype call for CYP2D6_10B / S1 --------------
$markerTable = $wb.Worksheets.Item("marker_table")
$markerTable.Range("G3").Value = "CT"
$markerTable.Range("H3").Value = "heterozygous"

# --- genotype_result: overall sample genotype --------------------------------
$genotypeResult = $wb.Worksheets.Item("genotype_result")
$genotypeResult.Range("B2").Value = "*1/*10B"
